# Update cryptos list (Price / Volume(1h) columns, plus two row-identity
# swaps) to match the latest scrape, mirroring the upstream GitHub Actions
# commit. Numeric-looking Price strings are written with a leading
# apostrophe so Excel keeps them as text (preserving formatting such as
# "9.20" or "3.86.xx" style grouped numbers) instead of silently coercing
# them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.351.76'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '3.648.26'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''195.74'
$ws.Range('E5').Value = '  +5.84%  '
$ws.Range('D6').Value = '''577.95'
$ws.Range('E6').Value = '  -2.02%  '
$ws.Range('D7').Value = '3.641.50'
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').Value = '''0.622'
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Value = '''0.682'
$ws.Range('E10').Value = '  +0.55%  '
$ws.Range('E11').Value = '  +4.91%  '
$ws.Range('D12').Value = '''57.39'
$ws.Range('E12').Value = '  +6.19%  '
$ws.Range('E13').Value = '  +15.77%  '
$ws.Range('D14').Value = '''10.19'
$ws.Range('E14').Value = '  +1.78%  '
$ws.Range('D15').Value = '4.219.84'
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('D16').Value = '3.631.13'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '''12.59'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('D19').Value = '68.199.19'
$ws.Range('E19').Value = '  +1.23%  '
$ws.Range('D20').Value = '''18.64'
$ws.Range('E20').Value = '  +0.78%  '
$ws.Range('D21').Value = '''1.09'
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('D22').Value = '''404.85'
$ws.Range('E22').Value = '  +1.39%  '
$ws.Range('E23').Value = '  +23.33%  '
$ws.Range('E24').Value = '  -1.97%  '
$ws.Range('D25').Value = '''86.28'
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('E26').Value = '  +3.11%  '
$ws.Range('D27').Value = '''12.67'
$ws.Range('E27').Value = '  +2.00%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '''3.86'
$ws.Range('E28').Value = '  +5.93%  '
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').Value = '''6.13'
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').Value = '''8.15'
$ws.Range('E30').Value = '  +19.89%  '
$ws.Range('D31').Value = '''9.20'
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('D32').Value = '''31.80'
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('D33').Value = '''693.44'
$ws.Range('E33').Value = '  +16.43%  '
$ws.Range('D34').Value = '''12.26'
$ws.Range('E34').Value = '  +2.56%  '
$ws.Range('D35').Value = '''0.119'
$ws.Range('E35').Value = '  +4.88%  '
$ws.Range('D36').Value = '''64.85'
$ws.Range('E36').Value = '  -3.16%  '
$ws.Range('D37').Value = '''42.96'
$ws.Range('E37').Value = '  +2.91%  '
$ws.Range('D38').Value = '''0.420'
$ws.Range('E38').Value = '  +10.99%  '
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('D40').Value = '0.0₃0799'
$ws.Range('E40').Value = '  +7.36%  '
$ws.Range('E41').Value = '  +18.98%  '
$ws.Range('E42').Value = '  +2.26%  '
$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D43').Value = '''3.13'
$ws.Range('E43').Value = '  +11.87%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '3.191.29'
$ws.Range('E44').Value = '  +17.58%  '
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').Value = '''2.93'
$ws.Range('E46').Value = '  +24.21%  '
$ws.Range('E47').Value = '  +2.28%  '
$ws.Range('E48').Value = '  +2.00%  '
$ws.Range('D49').Value = '''8.90'
$ws.Range('E49').Value = '  +7.44%  '
$ws.Range('D50').Value = '''3.14'
$ws.Range('E50').Value = '  +3.45%  '
$ws.Range('D51').Value = '''143.19'
$ws.Range('E51').Value = '  +3.65%  '
